$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EnvironmentProperties")
$ws.Range("B18").Value = "Aggressive(bool)"
$ws.Range("C18").Value = "AgroRange(SDL_Rect)"
Write-Host $ws.Range("B18").Value2
Write-Host $ws.Range("C18").Value2
